$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear AgTests (F) / AgPosit (G) columns for rows 393-483 (data reported as not-applicable for this range)
$ws.Range("F393:G483").ClearContents()

# Update AgTests (F) / AgPosit (G) values for rows 484-497
$ws.Cells.Item(484, 6).Value = 16238
$ws.Cells.Item(484, 7).Value = 22
$ws.Cells.Item(485, 6).Value = 27496
$ws.Cells.Item(485, 7).Value = 28
$ws.Cells.Item(486, 6).Value = 17678
$ws.Cells.Item(486, 7).Value = 14
$ws.Cells.Item(487, 6).Value = 13572
$ws.Cells.Item(487, 7).Value = 18
$ws.Cells.Item(488, 6).Value = 12548
$ws.Cells.Item(488, 7).Value = 16
$ws.Cells.Item(489, 6).Value = 23750
$ws.Cells.Item(489, 7).Value = 20
$ws.Cells.Item(490, 6).Value = 20292
$ws.Cells.Item(490, 7).Value = 30
$ws.Cells.Item(491, 6).Value = 18830
$ws.Cells.Item(491, 7).Value = 20
$ws.Cells.Item(492, 6).Value = 26080
$ws.Cells.Item(492, 7).Value = 22
$ws.Cells.Item(493, 6).Value = 15148
$ws.Cells.Item(493, 7).Value = 14
$ws.Cells.Item(494, 6).Value = 12092
$ws.Cells.Item(494, 7).Value = 12
$ws.Cells.Item(495, 6).Value = 19554
$ws.Cells.Item(495, 7).Value = 26
$ws.Cells.Item(496, 6).Value = 15114
$ws.Cells.Item(496, 7).Value = 28
$ws.Cells.Item(497, 6).Value = 13742
$ws.Cells.Item(497, 7).Value = 16

# Append new row 498
$ws.Cells.Item(498, 1).Value = 44392
$ws.Cells.Item(498, 2).Value = 392034
$ws.Cells.Item(498, 3).Value = 6155
$ws.Cells.Item(498, 4).Value = 34
$ws.Cells.Item(498, 5).Value = 12524
$ws.Cells.Item(498, 6).Value = 13924
$ws.Cells.Item(498, 7).Value = 12
